$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - update kickoff time and odds
$ws.Range("C10").Value = "23:07"

$ws.Range("I10").Value = 2.75
$ws.Range("J10").Value = 3.1
$ws.Range("L10").Value = 3.3
$ws.Range("N10").Value = 7.2
$ws.Range("P10").Value = 3.2
$ws.Range("R10").Value = 1.82
$ws.Range("T10").Value = 2.65
$ws.Range("U10").Value = 1.7
$ws.Range("V10").Value = 2.02
$ws.Range("W10").Value = 8.25
$ws.Range("X10").Value = 12
$ws.Range("AA10").Value = 20
$ws.Range("AB10").Value = 29
$ws.Range("AC10").Value = 7.2
$ws.Range("AE10").Value = 13
$ws.Range("AH10").Value = 9.25
$ws.Range("AI10").Value = 14.5
$ws.Range("AL10").Value = 22
$ws.Range("AM10").Value = 29
$ws.Range("AN10").Value = 4.4
$ws.Range("AO10").Value = 13.5
$ws.Range("AP10").Value = 22
$ws.Range("AQ10").Value = 60
$ws.Range("AR10").Value = 100
$ws.Range("AS10").Value = 300
$ws.Range("AT10").Value = 2.65
$ws.Range("AX10").Value = 15
$ws.Range("AY10").Value = 22
$ws.Range("AZ10").Value = 65
$ws.Range("BA10").Value = 100

# Row 11 - update kickoff time
$ws.Range("C11").Value = "23:11"
